# Update "想去人数" (F column) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 677
$ws1.Range("F6").Value  = 810
$ws1.Range("F10").Value = 1171
$ws1.Range("F12").Value = 357
$ws1.Range("F14").Value = 156
$ws1.Range("F15").Value = 104
$ws1.Range("F23").Value = 569

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value  = 207
$ws2.Range("F13").Value = 46

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value  = 677
$ws4.Range("F10").Value = 810
$ws4.Range("F14").Value = 1171
$ws4.Range("F18").Value = 357
$ws4.Range("F21").Value = 156
$ws4.Range("F22").Value = 104
$ws4.Range("F27").Value = 207
$ws4.Range("F32").Value = 46
$ws4.Range("F36").Value = 569
